# Diccionario_BiciAlpes.xlsx formatting touch-up
# ------------------------------------------------
# 1) The header row + data rows (A1:C19) used a theme-linked "black" font
#    color; switch them to an explicit black RGB color so the text color
#    no longer depends on the workbook theme.
# 2) A handful of row heights were nudged down slightly (header row, and
#    the empty trailing rows 20-24) to tighten up the layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Make the header/body font color an explicit black (was theme black) ---
$ws.Range("A1:C19").Font.Color = 0   # RGB(0,0,0) -> serialized as rgb="FF000000"

# --- 2) Row height tweaks ---
$ws.Rows.Item(1).RowHeight = 20.25

$ws.Rows.Item(20).RowHeight = 18.75
$ws.Rows.Item(21).RowHeight = 18.75
$ws.Rows.Item(22).RowHeight = 18.75
$ws.Rows.Item(23).RowHeight = 18.75
$ws.Rows.Item(24).RowHeight = 18.75
